$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -11.96559999999999

$ws.Range("A9").Value = -22.04460000000001
$ws.Range("C9").Value = -12.43739999999999
$ws.Range("D9").Value = -8.632200000000006

$ws.Range("A18").Value = -22.15330000000001

$ws.Range("A20").Value = -20.15999999999999

$ws.Range("C23").Value = -12.2516

$ws.Range("C24").Value = -13.2469

$ws.Range("C26").Value = -12.3405

$ws.Range("A27").Value = -22.00389999999999

$ws.Range("D32").Value = -7.158499999999997

$ws.Range("C34").Value = -11.69290000000001

$ws.Range("C35").Value = -12.3555

$ws.Range("D38").Value = -7.432799999999999

$ws.Range("D45").Value = -7.166099999999996

$ws.Range("C48").Value = -11.4456

$ws.Range("D51").Value = -8.298400000000004

$ws.Range("C52").Value = -11.0746

$ws.Range("D57").Value = -7.745999999999996

$ws.Range("D64").Value = -7.081599999999995

$ws.Range("C66").Value = -11.0149

$ws.Range("C67").Value = -11.4272

$ws.Range("A69").Value = -21.47549999999998

$ws.Range("A76").Value = -19.65049999999999

$ws.Range("C80").Value = -13.457

$ws.Range("A82").Value = -21.71160000000001

$ws.Range("D93").Value = -6.920099999999993

$ws.Range("C99").Value = -12.3183
